$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 57903.383
$ws.Range("I33").Value = 75209.5
$ws.Range("J33").Value = 216.33333
$ws.Range("K33").Value = 75209.5
$ws.Range("L33").Value = 216.33333
$ws.Range("M33").Value = -74980.5
$ws.Range("N33").Value = -674.3333299999999
$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 10000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30584
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31068
$ws.Range("H76").Value = 5627674
$ws.Range("I76").Value = 4845.3335
$ws.Range("K76").Value = 4845.3335
$ws.Range("M76").Value = -4530.3335
$ws.Range("H79").Value = 5627674
$ws.Range("I79").Value = 4845.3335
$ws.Range("K79").Value = 4845.3335
$ws.Range("M79").Value = -3753.3335
$ws.Range("H106").Value = 80467.46000000001
$ws.Range("I106").Value = 3839.75
$ws.Range("K106").Value = 3839.75
$ws.Range("M106").Value = -3208.75
$ws.Range("H113").Value = 71431580
$ws.Range("I113").Value = 111113520
$ws.Range("J113").Value = 4100.2
$ws.Range("K113").Value = 111113520
$ws.Range("L113").Value = 4100.2
$ws.Range("M113").Value = -111110266
$ws.Range("N113").Value = -10608.2
$ws.Range("H137").Value = 4704.125
$ws.Range("I137").Value = 4416.4165
$ws.Range("J137").Value = 4991.8335
$ws.Range("K137").Value = 13249.2495
$ws.Range("L137").Value = 14975.5005
$ws.Range("M137").Value = -10699.2495
$ws.Range("N137").Value = -20075.5005
$ws.Range("H138").Value = 6171.5293
$ws.Range("I138").Value = 5434.5
$ws.Range("J138").Value = 6687.45
$ws.Range("K138").Value = 16303.5
$ws.Range("L138").Value = 20062.35
$ws.Range("M138").Value = -11163.5
$ws.Range("N138").Value = -30342.35
$ws.Range("H141").Value = 6079.857
$ws.Range("I141").Value = 6027.2964
$ws.Range("J141").Value = 7499
$ws.Range("K141").Value = 18081.8892
$ws.Range("L141").Value = 22497
$ws.Range("M141").Value = -12901.8892
$ws.Range("N141").Value = -32857

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31899.047
$ws.Range("I32").Value = 21747.04
$ws.Range("J32").Value = 48396.062
$ws.Range("K32").Value = 21747.04
$ws.Range("L32").Value = 48396.062
$ws.Range("M32").Value = -21460.04
$ws.Range("N32").Value = -48970.062
$ws.Range("H132").Value = 478629.88
$ws.Range("J132").Value = 1432356.9
$ws.Range("L132").Value = 4297070.699999999
$ws.Range("N132").Value = -4302130.699999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29005.75
$ws.Range("I20").Value = 5341.3335
$ws.Range("K20").Value = 5341.3335
$ws.Range("M20").Value = -5094.3335
$ws.Range("H132").Value = 52500
$ws.Range("J132").Value = 52500
$ws.Range("L132").Value = 52500
$ws.Range("N132").Value = -62620
$ws.Range("H134").Value = 3678.2964
$ws.Range("I134").Value = 3331.261
$ws.Range("K134").Value = 9993.782999999999
$ws.Range("M134").Value = -7458.782999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4989.6
$ws.Range("I31").Value = 4499.5
$ws.Range("J31").Value = 5167.8184
$ws.Range("K31").Value = 4499.5
$ws.Range("L31").Value = 5167.8184
$ws.Range("M31").Value = -4204.5
$ws.Range("N31").Value = -5757.8184
$ws.Range("H34").Value = 4989.6
$ws.Range("I34").Value = 4499.5
$ws.Range("J34").Value = 5167.8184
$ws.Range("K34").Value = 4499.5
$ws.Range("L34").Value = 5167.8184
$ws.Range("M34").Value = -4297.5
$ws.Range("N34").Value = -5571.8184
$ws.Range("H134").Value = 3773.8948
$ws.Range("J134").Value = 7013.2856
$ws.Range("L134").Value = 21039.8568
$ws.Range("N134").Value = -26109.8568
$ws.Range("H141").Value = 169066.12
$ws.Range("J141").Value = 169066.12
$ws.Range("L141").Value = 169066.12
$ws.Range("N141").Value = -179426.12

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4259.5
$ws.Range("I3").Value = 4259.5
$ws.Range("K3").Value = 12778.5
$ws.Range("M3").Value = -12666.5
$ws.Range("H41").Value = 51.6
$ws.Range("I41").Value = 51.6
$ws.Range("K41").Value = 154.8
$ws.Range("M41").Value = 183.2
$ws.Range("H113").Value = 2566376.5
$ws.Range("I113").Value = 1349.5
$ws.Range("J113").Value = 3207633.5
$ws.Range("K113").Value = 4048.5
$ws.Range("L113").Value = 9622900.5
$ws.Range("M113").Value = -1878.5
$ws.Range("N113").Value = -9627240.5
$ws.Range("H129").Value = 1434.7273
$ws.Range("J129").Value = 2993
$ws.Range("L129").Value = 8979
$ws.Range("N129").Value = -18979
$ws.Range("H132").Value = 1886.5333
$ws.Range("I132").Value = 1650.25
$ws.Range("J132").Value = 2156.5715
$ws.Range("K132").Value = 14852.25
$ws.Range("L132").Value = 19409.1435
$ws.Range("M132").Value = -12322.25
$ws.Range("N132").Value = -24469.1435

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3456.4167
$ws.Range("I80").Value = 2997.8
$ws.Range("J80").Value = 3784
$ws.Range("K80").Value = 2997.8
$ws.Range("L80").Value = 3784
$ws.Range("M80").Value = -1999.8
$ws.Range("N80").Value = -5780
$ws.Range("H83").Value = 3456.4167
$ws.Range("I83").Value = 2997.8
$ws.Range("J83").Value = 3784
$ws.Range("K83").Value = 14989
$ws.Range("L83").Value = 18920
$ws.Range("M83").Value = -9997
$ws.Range("N83").Value = -28904
$ws.Range("H132").Value = 423760.72
$ws.Range("I132").Value = 594074.3
$ws.Range("J132").Value = 10142
$ws.Range("K132").Value = 1782222.9
$ws.Range("L132").Value = 30426
$ws.Range("M132").Value = -1779692.9
$ws.Range("N132").Value = -35486

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 7999
$ws.Range("I18").Value = 7999
$ws.Range("K18").Value = 7999
$ws.Range("M18").Value = -7827
$ws.Range("H61").Value = 417912.9
$ws.Range("I61").Value = 501157.1
$ws.Range("J61").Value = 1692
$ws.Range("K61").Value = 501157.1
$ws.Range("L61").Value = 1692
$ws.Range("M61").Value = -500955.1
$ws.Range("N61").Value = -2096
$ws.Range("H113").Value = 417912.9
$ws.Range("I113").Value = 501157.1
$ws.Range("J113").Value = 1692
$ws.Range("K113").Value = 501157.1
$ws.Range("L113").Value = 1692
$ws.Range("M113").Value = -498987.1
$ws.Range("N113").Value = -6032
$ws.Range("H132").Value = 159074.19
$ws.Range("I132").Value = 266818.4
$ws.Range("J132").Value = 7434.185
$ws.Range("K132").Value = 800455.2000000001
$ws.Range("L132").Value = 22302.555
$ws.Range("M132").Value = -797925.2000000001
$ws.Range("N132").Value = -27362.555
$ws.Range("H136").Value = 25648666
$ws.Range("I136").Value = 43486836
$ws.Range("K136").Value = 130460508
$ws.Range("M136").Value = -130457958

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 483492.94
$ws.Range("I132").Value = 723025.9
$ws.Range("J132").Value = 4427.143
$ws.Range("K132").Value = 2169077.7
$ws.Range("L132").Value = 13281.429
$ws.Range("M132").Value = -2166547.7
$ws.Range("N132").Value = -18341.429
$ws.Range("H136").Value = 6218.4736
$ws.Range("I136").Value = 7076.8667
$ws.Range("K136").Value = 21230.6001
$ws.Range("M136").Value = -18680.6001
